# Weekly refresh of "Ramas de apio" price records.
# The underlying data rows (2-19) get re-mapped to a new set of
# positions (a cyclic re-shuffle produced by the upstream weekly
# consolidation job). We snapshot every row first (so we never read a
# value that has already been overwritten), then write each snapshot
# back out to its new row according to the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 19
$lastCol = 18

# new row -> source row (i.e. new row R gets the full contents that used
# to live in row Mapping[R])
$mapping = @{
    2  = 5
    3  = 7
    4  = 4
    5  = 18
    6  = 11
    7  = 13
    8  = 3
    9  = 15
    10 = 6
    11 = 9
    12 = 16
    13 = 8
    14 = 17
    15 = 19
    16 = 10
    17 = 12
    18 = 14
    19 = 2
}

# Snapshot all existing rows before we start writing anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $vals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $vals
}

# Write the snapshotted rows back out in their new positions.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
